$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01886266666666667
$ws.Range("H2").Value = 0.056588
$ws.Range("M2").Value = 0.299213
$ws.Range("N2").Value = 0.598426
$ws.Range("O2").Value = 0.09740004925129409
$ws.Range("P2").Value = 0.07136653219636822
$ws.Range("Q2").Value = 0.005643955081333333
$ws.Range("R2").Value = 0.033863730488
$ws.Range("S2").Value = 0.09740004925129409
$ws.Range("T2").Value = 0.07136653219636822

# Row 3
$ws.Range("G3").Value = 0.01886266666666667
$ws.Range("H3").Value = 0.056588
$ws.Range("O3").Value = 0.1329300998920193
$ws.Range("P3").Value = 0.1460999300318762
$ws.Range("Q3").Value = 0.007702783710222223
$ws.Range("R3").Value = 0.069325053392
$ws.Range("S3").Value = 0.1329300998920193
$ws.Range("T3").Value = 0.1460999300318762

# Row 4
$ws.Range("G4").Value = 0.01886266666666667
$ws.Range("H4").Value = 0.056588
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04598133333333334
$ws.Range("N4").Value = 0.137944
$ws.Range("O4").Value = 0.01496787950826614
$ws.Range("P4").Value = 0.01645079745414774
$ws.Range("Q4").Value = 0.0008673305635555557
$ws.Range("R4").Value = 0.007805975072
$ws.Range("S4").Value = 0.01496787950826614
$ws.Range("T4").Value = 0.01645079745414774

# Row 5
$ws.Range("G5").Value = 0.01886266666666667
$ws.Range("H5").Value = 0.056588
$ws.Range("M5").Value = 0.5315415
$ws.Range("N5").Value = 1.063083
$ws.Range("O5").Value = 0.1730278038691726
$ws.Range("P5").Value = 0.1267801652115913
$ws.Range("Q5").Value = 0.010026290134
$ws.Range("R5").Value = 0.060157740804
$ws.Range("S5").Value = 0.1730278038691726
$ws.Range("T5").Value = 0.1267801652115913

# Row 6
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.01886266666666667
$ws.Range("H6").Value = 0.056588
$ws.Range("M6").Value = 1.667958333333333
$ws.Range("N6").Value = 5.003875
$ws.Range("O6").Value = 0.5429550982603464
$ws.Range("P6").Value = 0.5967474780408973
$ws.Range("Q6").Value = 0.03146214205555555
$ws.Range("R6").Value = 0.2831592785
$ws.Range("S6").Value = 0.5429550982603464
$ws.Range("T6").Value = 0.5967474780408973

# Row 7
$ws.Range("D7").Value = "Neutrophils"
$ws.Range("G7").Value = 0.01886266666666667
$ws.Range("H7").Value = 0.056588
$ws.Range("M7").Value = 0.118945
$ws.Range("N7").Value = 0.356835
$ws.Range("O7").Value = 0.0387190692189015
$ws.Range("P7").Value = 0.04255509706511925
$ws.Range("Q7").Value = 0.002243619886666667
$ws.Range("R7").Value = 0.02019257898
$ws.Range("S7").Value = 0.0387190692189015
$ws.Range("T7").Value = 0.04255509706511925

